$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values for rows 2-10 (B:G), computed from new forecast-error data
$data = @(
    @(0.03678685664255971,  0.6869725908349552,  1.143126198255408,   1.069170799383994,   1.079170229362562,   51),
    @(-0.183591446892805,   0.4867160004773056,  0.4160893806761118,  0.6450499055701906,  0.6246497478211112,  50),
    @(-0.03936449295927064, 0.5795287071841984,  0.7040057845192284,  0.839050525605716,   0.8468120959350416,  49),
    @(-0.1295444554373254,  0.4050807838040723,  0.2941531664628659,  0.5423588908304776,  0.5322338956591262,  48),
    @(0.00637098485360963,  0.5886252521369484,  0.7430823678720198,  0.8620222548588985,  0.8713178872051619,  47),
    @(-0.132221774258578,   0.3680198856531414,  0.2162217108530362,  0.4649964632693846,  0.4507277870335607,  46),
    @(0.07610304665680741,  0.5642969124478534,  0.6783264378029779,  0.8236057538670901,  0.8293489284340051,  45),
    @(-0.07952857203228597, 0.3289968720127114,  0.1737617843583987,  0.4168474353506312,  0.4139213380474263,  44),
    @(0.1447151924720829,   0.5826609323638429,  0.8037054263539364,  0.8964961942774416,  0.8952095577377838,  43)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

# New row 11: Q9 (copy formatting/style from A10, then overwrite the value)
$ws.Cells.Item(10, 1).Copy($ws.Cells.Item(11, 1))
$ws.Cells.Item(11, 1).Value = "Q9"
$ws.Cells.Item(11, 2).Value = -0.06616833373077866
$ws.Cells.Item(11, 3).Value = 0.3497216611734643
$ws.Cells.Item(11, 4).Value = 0.2089885319410827
$ws.Cells.Item(11, 5).Value = 0.4571526352774122
$ws.Cells.Item(11, 6).Value = 0.4578217756651212
$ws.Cells.Item(11, 7).Value = 42
